$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Six existing rows get a new "FL(5-0-0)" leave entry (5 days) recorded in
#    the PARTICULARS / Absence-Undertime-W/Pay columns.
# ---------------------------------------------------------------------------
$flRows = 22,35,48,61,74,87
foreach ($r in $flRows) {
    $ws.Range("B$r").Value = "FL(5-0-0)"
    $ws.Range("D$r").Value = 5
}

# ---------------------------------------------------------------------------
# 2) Insert a brand-new row at worksheet row 541 (this pushes the existing
#    rows 541-604 down to 542-605, and the table grows from A8:K604 to
#    A8:K605). The new row records a "FL(2-0-0)" leave entry (2 days).
# ---------------------------------------------------------------------------
$ws.Range("A541").EntireRow.Insert()

# Copy the (now shifted) row below it to pick up identical cell formatting
# for the freshly inserted, still-unformatted row.
$ws.Range("A543:K543").Copy()
$ws.Range("A541:K541").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Re-attach the new row to Table1 so the table range grows to A8:K605.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A8:K605"))

# Restore the calculated "EARNED " helper-column formula on the new row and
# on the (re-numbered) final table row, whose structured reference can go
# stale immediately after the insert/resize.
$ws.Range("G541").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G605").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Fill in the new leave entry itself.
$ws.Range("B541").Value = "FL(2-0-0)"
$ws.Range("D541").Value = 2

# ---------------------------------------------------------------------------
# 3) Recalculate so dependent totals (e.g. the summary BALANCE formula in
#    E9) reflect the newly-added absence/undertime values.
# ---------------------------------------------------------------------------
$excel.CalculateFull()

# ---------------------------------------------------------------------------
# 4) Restore on-screen selection state to match the latest save.
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.Panes.Item(2).Activate()
$ws.Range("C587").Select()
$win.Panes.Item(1).Activate()
$ws.Range("E2:E3").Select()
